$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$shp = $ws.Shapes.Item("Picture 3")
$shp.Height = $shp.Height - 3
Write-Host ("Height=" + $shp.Height)
